$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# "Changed layout - Rearranged components in BOM"
#
# Rows 34 and 35 swap places: the U1 / STM32G0B1CEU6 (MCU) line moves up to
# row 34 (gaining the taller 29pt row height its wrapped description needs)
# and the U2 / T2035H-6G (Thyristor) line drops down to row 35 (back to the
# sheet's default row height). The hyperlink attached to column E follows
# each part's row.
# ---------------------------------------------------------------------------

# Snapshot the current values of row 34 (U2 / Thyristor) ...
$a34 = $ws.Range("A34").Value()
$b34 = $ws.Range("B34").Value()
$c34 = $ws.Range("C34").Value()
$d34 = $ws.Range("D34").Value()
$e34 = $ws.Range("E34").Value()
$f34 = $ws.Range("F34").Value()

# ... and row 35 (U1 / MCU).
$a35 = $ws.Range("A35").Value()
$b35 = $ws.Range("B35").Value()
$c35 = $ws.Range("C35").Value()
$d35 = $ws.Range("D35").Value()
$e35 = $ws.Range("E35").Value()
$f35 = $ws.Range("F35").Value()

# Write the former row 35 (U1 / MCU) content into row 34.
$ws.Range("A34").Value = $a35
$ws.Range("B34").Value = $b35
$ws.Range("C34").Value = $c35
$ws.Range("D34").Value = $d35
$ws.Range("E34").Value = $e35
$ws.Range("F34").Value = $f35

# Write the former row 34 (U2 / Thyristor) content into row 35.
$ws.Range("A35").Value = $a34
$ws.Range("B35").Value = $b34
$ws.Range("C35").Value = $c34
$ws.Range("D35").Value = $d34
$ws.Range("E35").Value = $e34
$ws.Range("F35").Value = $f34

# Row heights follow the content: row 34 now holds the long wrapped
# description (29pt tall), row 35 goes back to the sheet's default height.
$ws.Rows(34).RowHeight = 29
$ws.Rows(35).AutoFit()

# Re-point the hyperlinks living on column E so they stay attached to the
# correct part number after the swap.
foreach ($hl in $ws.Hyperlinks) {
    $addr = $hl.Range().Address()
    if ($addr -eq '$E$34') {
        $hl.Address = "https://octopart.com/stm32g0b1ceu6-stmicroelectronics-116363364?r=sp"
    } elseif ($addr -eq '$E$35') {
        $hl.Address = "https://octopart.com/t2035h-6g-stmicroelectronics-9417760"
    }
}

# Update the view state left in the sheet: scrolled further down with B32
# selected.
$ws.Range("B32").Select()
